$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the newly added "NewText" column (C)
$ws.Range("C1").Value = "NewText"

# Replace existing text samples (column A) with new test content
$ws.Range("A2").Value = "Trust the magic of new beginnings!!!"
$ws.Range("A3").Value = "Keep going..!!!"

# New row for the added test case (Notepad_TC03_replaceText)
$ws.Range("A4").Value = "Good Things take time!!!!"
$ws.Range("B4").Value = "NotepadFile4"
$ws.Range("C4").Value = "You are enough just as you are!!!"

# Match the cursor/selection position left behind by the author's save
$ws.Range("A9").Select() | Out-Null
